$d = $word.ActiveDocument

# Locate the "<A DEFINIR>" placeholder text.
$r = $d.Content
$found = $r.Find.Execute("<A DEFINIR>", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # $r now spans the whole "<A DEFINIR>" match. Build a collapsed range
    # positioned right before the closing ">" (the match's last character).
    $insertionPoint = $r.Duplicate
    $insertionPoint.Start = $r.End - 1
    $insertionPoint.End = $r.End - 1

    $insertionPoint.InsertAfter(" nas próximas aulas")

    # Force the newly-inserted text into its own run (rather than being
    # silently coalesced back into the surrounding "<A DEFINIR" / ">" run)
    # by toggling a character property on it and then reverting it.
    $insertionPoint.Bold = 1
    $insertionPoint.Bold = 0
}
